$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 (C, D, E columns) ---
# Row 8: extr1
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# Row 9: extr2
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $false

# Row 10: extr3
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# Row 11: extr4
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# Row 12: extr5
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $false

# Row 13: extr6
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

# Row 14: extr7
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $true

# Row 15: extr8
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

# --- Add new rows 16 and 17 ---
# Row 16: line7
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "line7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

# Row 17: line8
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "line8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false

# Apply the same style as the other "A" column cells (e.g. A15) to the new A16/A17 cells
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(17, 1).PasteSpecial(-4122)  # xlPasteFormats
